$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Label" in H1 (copy the header style from G1 first)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"
$excel.CutCopyMode = 0

# New "Label" column values: 0 for control rows, 1 for MDD rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 1

# Refit prediction/error values (refitting NCDEs to individual patients)
$ws.Range("D3").Value = 0.4819464224482482
$ws.Range("E3").Value = 0.4819464224482482

$ws.Range("D4").Value = 0.4084233238187719
$ws.Range("E4").Value = 0.4084233238187719

$ws.Range("D7").Value = 0.7504624208557402
$ws.Range("E7").Value = 0.2495375791442598
